# Daily attendance processing - 2025-11-23 13:39:07
# Swap the order of the two comma-separated values in column "Recorded By" (G)
# for rows where the value currently starts with "dnasr281@gmail.com, ".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("dnasr281@gmail.com, ")) {
        $parts = $val.ToString().Split(",")
        if ($parts.Length -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            $cell.Value = "$second, $first"
        }
    }
}
